$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to retain text formatting so numeric-looking values
# (e.g. "1.000") are not auto-converted to numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '30.272.72'
$ws.Range("E2").Value = '  -0.14%  '

$ws.Range("D3").Value = '1.862.23'

$ws.Range("E4").Value = '  -0.02%  '

$ws.Range("D5").Value = '236.33'
$ws.Range("E5").Value = '  +0.44%  '

$ws.Range("D6").Value = '1.000'
$ws.Range("E6").Value = '  -0.03%  '

$ws.Range("D7").Value = '0.4706'
$ws.Range("E7").Value = '  +0.77%  '

$ws.Range("D8").Value = '0.2908'
$ws.Range("E8").Value = '  +2.35%  '

$ws.Range("D9").Value = '0.06535'
$ws.Range("E9").Value = '  -0.43%  '

$ws.Range("D10").Value = '21.81'
$ws.Range("E10").Value = '  +2.61%  '

$ws.Range("D11").Value = '0.07932'
$ws.Range("E11").Value = '  +0.64%  '

$ws.Range("D12").Value = '97.74'
$ws.Range("E12").Value = '  -0.56%  '

$ws.Range("D13").Value = '1.860.53'
$ws.Range("E13").Value = '  -0.56%  '

$ws.Range("D14").Value = '5.148'
$ws.Range("E14").Value = '  +0.53%  '

$ws.Range("D15").Value = '0.6811'
$ws.Range("E15").Value = '  +0.75%  '

$ws.Range("D16").Value = '263.90'
$ws.Range("E16").Value = '  -6.23%  '

$ws.Range("D17").Value = '30.256.12'
$ws.Range("E17").Value = '  -0.15%  '

$ws.Range("D18").Value = '13.77'
$ws.Range("E18").Value = '  +8.46%  '

$ws.Range("D19").Value = '1.000'
$ws.Range("E19").Value = '  -0.07%  '

$ws.Range("D20").Value = '0.000007448'
$ws.Range("E20").Value = '  +2.09%  '

$ws.Range("D21").Value = '2.106.72'
$ws.Range("E21").Value = '  -0.15%  '

$ws.Range("D22").Value = '1.001'
$ws.Range("E22").Value = '  +0.04%  '

$ws.Range("D23").Value = '5.273'
$ws.Range("E23").Value = '  -4.26%  '

$ws.Range("D24").Value = '6.170'

$ws.Range("D25").Value = '167.39'
$ws.Range("E25").Value = '  +1.53%  '

$ws.Range("D26").Value = '9.201'
$ws.Range("E26").Value = '  -0.44%  '

$ws.Range("E27").Value = '  -1.54%  '

$ws.Range("D28").Value = '1.953'
$ws.Range("E28").Value = '  +1.08%  '

$ws.Range("E29").Value = '  +1.52%  '

$ws.Range("D30").Value = '0.09853'
$ws.Range("E30").Value = '  +1.34%  '

$ws.Range("D31").Value = '4.354'
$ws.Range("E31").Value = '  -1.62%  '

$ws.Range("E32").Value = '  -0.48%  '

$ws.Range("D33").Value = '4.034'
$ws.Range("E33").Value = '  -2.00%  '

$ws.Range("D34").Value = '0.04716'
$ws.Range("E34").Value = '  +0.52%  '

$ws.Range("D35").Value = '1.128'
$ws.Range("E35").Value = '  +0.10%  '

$ws.Range("D36").Value = '0.6994'
$ws.Range("E36").Value = '  -0.99%  '

$ws.Range("D37").Value = '2.706'
$ws.Range("E37").Value = '  -0.24%  '

$ws.Range("D38").Value = '0.01876'
$ws.Range("E38").Value = '  +0.77%  '

$ws.Range("D39").Value = '2.621'
$ws.Range("E39").Value = '  +3.12%  '

$ws.Range("D40").Value = '6.336'
$ws.Range("E40").Value = '  +0.69%  '

$ws.Range("D41").Value = '73.95'
$ws.Range("E41").Value = '  +0.62%  '

$ws.Range("D42").Value = '1.946'
$ws.Range("E42").Value = '  -0.43%  '

$ws.Range("D43").Value = '0.8453'
$ws.Range("E43").Value = '  -0.11%  '

$ws.Range("B44").Value = 'PaxDollar'
$ws.Range("C44").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D44").Value = '0.9995'
$ws.Range("E44").Value = '  -0.13%  '

$ws.Range("B45").Value = 'TheSandbox'
$ws.Range("C45").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D45").Value = '0.4157'
$ws.Range("E45").Value = '  -0.53%  '

$ws.Range("D46").Value = '103.23'
$ws.Range("E46").Value = '  -0.74%  '

$ws.Range("D47").Value = '7.158'
$ws.Range("E47").Value = '  -0.55%  '

$ws.Range("D48").Value = '943.90'
$ws.Range("E48").Value = '  +1.12%  '

$ws.Range("D49").Value = '9.222'
$ws.Range("E49").Value = '  +0.94%  '

$ws.Range("D50").Value = '34.16'
$ws.Range("E50").Value = '  +0.22%  '

$ws.Range("D51").Value = '0.05661'
$ws.Range("E51").Value = '  +0.54%  '
